$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the latest cryptocurrency price/volume snapshot.
# Price-column (D) values are written with an explicit Text number
# format so values like "1.001", "0.5700", or "27.908.35" are stored
# as literal text (matching the scraped data feed) rather than being
# auto-coerced into floating point numbers by Excel.

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "27.908.35"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -2.47%  "

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "1.793.49"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -0.60%  "

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.11%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "316.93"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.20%  "

$ws.Cells.Item(6, 5).Value = "  +0.07%  "

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.5312"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -2.59%  "

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.3931"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +3.57%  "

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.07447"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -0.86%  "

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "41.44"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -2.22%  "

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "1.084"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -2.73%  "

$ws.Cells.Item(12, 5).Value = "  +0.00%  "

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "6.176"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +0.33%  "

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "7.479"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +1.27%  "

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "20.33"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -1.69%  "

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "1.795.46"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -0.40%  "

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "88.27"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -2.14%  "

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "0.00001058"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -0.77%  "

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "0.06580"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +1.48%  "

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "17.17"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -0.79%  "

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "5.944"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.20%  "

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "27.951.36"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -2.40%  "

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "11.07"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -0.28%  "

$ws.Cells.Item(25, 5).Value = "  -0.33%  "

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "156.88"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -2.38%  "

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "20.14"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -1.46%  "

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "2.007.72"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -0.03%  "

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "2.285"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -3.13%  "

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "121.72"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -1.11%  "

$ws.Cells.Item(31, 2).Value = "Stellar"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "0.1085"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +2.72%  "

$ws.Cells.Item(32, 2).Value = "ImmutableX"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "1.093"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -1.36%  "

$ws.Cells.Item(33, 5).Value = "  -0.24%  "

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "5.490"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -2.57%  "

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "0.07109"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +7.17%  "

$ws.Cells.Item(36, 5).Value = "  -2.47%  "

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "5.096"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +1.53%  "

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.02269"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -1.39%  "

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "8.365"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -4.51%  "

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "11.21"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -0.50%  "

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "0.6106"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -2.21%  "

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "1.177"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -1.43%  "

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "1.418"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -1.32%  "

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "13.28"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.11%  "

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "3.679"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -0.46%  "

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "0.5700"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -2.61%  "

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "125.08"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -1.26%  "

$ws.Cells.Item(48, 5).Value = "  +1.63%  "

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "1.913"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -1.60%  "

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "0.06805"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -1.27%  "

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "71.20"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -1.49%  "

